# Applies the OOXML changes described by the commit:
#  - Merge several list-item paragraphs in the body that were split into
#    two runs (text + trailing ".") into a single run each.
#  - Merge the "...cuando se" + "a" + " finalizado" + "." runs into one.
#  - Update the "Salida a producción 19 -01 - 2022" comment to "...2023".

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute(
        $oldText,   # FindText
        $false,     # MatchCase
        $false,     # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        $newText,   # ReplaceWith
        2           # Replace (wdReplaceAll)
    ) | Out-Null
}

# --- Body paragraphs: collapse split runs ("text" + ".") into one run ---

Replace-Text "EL alumno podrá ver los cursos en los cuales se puede matricular." `
             "EL alumno podrá ver los cursos en los cuales se puede matricular."

Replace-Text "El alumno puede solicitar ser matriculado a un curso, una vez aprobada la matricula podrá dar inicio al curso solicitado." `
             "El alumno puede solicitar ser matriculado a un curso, una vez aprobada la matricula podrá dar inicio al curso solicitado."

Replace-Text "Como alumno externo debo poder acceder al curso adquirido." `
             "Como alumno externo debo poder acceder al curso adquirido."

Replace-Text "Si cuento con un plan carreara como alumno interno debo visualizar el contenido de mi pénsum." `
             "Si cuento con un plan carreara como alumno interno debo visualizar el contenido de mi pénsum."

Replace-Text "Cada tarjeta debe iniciar con una barra de color verde en la parte superior." `
             "Cada tarjeta debe iniciar con una barra de color verde en la parte superior."

Replace-Text "Cada tarjeta debe tener un numero dentro de un cuadro gris en la parte superior izquierda." `
             "Cada tarjeta debe tener un numero dentro de un cuadro gris en la parte superior izquierda."

Replace-Text "Debajo del número debe tener un título del tema que se abordara." `
             "Debajo del número debe tener un título del tema que se abordara."

Replace-Text "Al hacer clic en una tarjeta esta despliega el contenido de forma descendente en una modal." `
             "Al hacer clic en una tarjeta esta despliega el contenido de forma descendente en una modal."

Replace-Text "En el costado superior derecha contendrá una X para cerrar la modal o haciendo clic en otra tarjeta o perdiendo el foco." `
             "En el costado superior derecha contendrá una X para cerrar la modal o haciendo clic en otra tarjeta o perdiendo el foco."

Replace-Text "El contenido de la modal será; El número correspondiente dentro de un cuadro gris justificado a la izquierda, a su derecha tendrá el título del contenido con una barra de color verde debajo." `
             "El contenido de la modal será; El número correspondiente dentro de un cuadro gris justificado a la izquierda, a su derecha tendrá el título del contenido con una barra de color verde debajo."

# This one keeps the "check" spell-check run (spellStart/spellEnd) intact;
# only merges the trailing " de color gris...a finalizado." runs.
Replace-Text " de color gris y cambiará a corlo verde cuando sea finalizado." `
             " de color gris y cambiará a corlo verde cuando sea finalizado."

# --- Comment: bump the year in the "Salida a producción" note ---

for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $c = $d.Comments.Item($i)
    if ($c.Range.Text -eq "Salida a producción 19 -01 - 2022") {
        $c.Range.Text = "Salida a producción 19 -01 - 2023"
        break
    }
}
